$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume table with the latest scraped values.
# Column D ("Price") holds numeric-looking text (e.g. "1.000", "72.00",
# "28.154.58"); assigning such strings straight to .Value would let Excel
# reinterpret them as numbers and silently drop significant trailing/
# thousands-style zeros, so each Price cell is forced to Text format
# immediately before its new value is written.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.154.58"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.804.86"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.69"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5483"
$ws.Range("E7").Value = "  +3.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3796"
$ws.Range("E8").Value = "  +0.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07479"
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.16"
$ws.Range("E10").Value = "  -0.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.096"
$ws.Range("E11").Value = "  -2.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.206"
$ws.Range("E13").Value = "  +0.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.56"
$ws.Range("E14").Value = "  -2.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.360"
$ws.Range("E15").Value = "  -2.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.804.11"
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.03"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06525"
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.47"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.935"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.185.11"
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.093"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.33"
$ws.Range("E26").Value = "  -2.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.50"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.012.47"
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "122.18"
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1129"
$ws.Range("E31").Value = "  +9.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.124"
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.659"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.586"
$ws.Range("E34").Value = "  -1.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06940"
$ws.Range("E35").Value = "  +6.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2234"
$ws.Range("E36").Value = "  -2.32%  "
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.096"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.486"
$ws.Range("E39").Value = "  -4.21%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.19"
$ws.Range("E40").Value = "  -2.30%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6187"
$ws.Range("E41").Value = "  -1.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.177"
$ws.Range("E42").Value = "  -2.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.421"
$ws.Range("E43").Value = "  +1.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.40"
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5768"
$ws.Range("E46").Value = "  -2.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "124.98"
$ws.Range("E47").Value = "  -1.06%  "
$ws.Range("E48").Value = "  +2.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.926"
$ws.Range("E49").Value = "  -2.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06822"
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.00"
$ws.Range("E51").Value = "  -1.54%  "
